# Adds a block of new "Math" unit tests (columns E/F) to Sheet1, plus a
# couple of extra tests in columns C/D (POW ^ / POW), matching the
# "Added more unit tests for basic functions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = -4152

# ---------------------------------------------------------------------
# Header row (E1:F1) -- same text/bold style as the existing B1/D1 "Value"
# header pair.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "Name"
$ws.Range("E1").Font.Bold = $true

$ws.Range("F1").Value = "Value"
$ws.Range("F1").Font.Bold = $true

# ---------------------------------------------------------------------
# Extra POW tests added to columns C/D (entered first, so the shared
# strings table lists "POW ^"/"POW" ahead of the E/F column labels).
# ---------------------------------------------------------------------

# Row 20: POW ^
$ws.Range("C20").Value = "POW ^"
$ws.Range("C20").HorizontalAlignment = $xlRight
$ws.Range("D20").Formula = "=2^3"

# Row 21: POW
$ws.Range("C21").Value = "POW"
$ws.Range("C21").HorizontalAlignment = $xlRight
$ws.Range("D21").Formula = "=POWER(2, 3)"

# ---------------------------------------------------------------------
# New column E labels + column F formulas (rows 2-11), each label cell
# right-aligned like the existing A/C label columns.
# ---------------------------------------------------------------------

# Row 2: TO DEG
$ws.Range("E2").Value = "TO DEG"
$ws.Range("E2").HorizontalAlignment = $xlRight
$ws.Range("F2").Formula = "=DEGREES(PI()/2)"

# Row 3: TO RAD
$ws.Range("E3").Value = "TO RAD"
$ws.Range("E3").HorizontalAlignment = $xlRight
$ws.Range("F3").Formula = "=RADIANS(90)"

# Row 4: blank label cell, kept right-aligned like its neighbours
$ws.Range("E4").HorizontalAlignment = $xlRight

# Row 5: QUOTIENT
$ws.Range("E5").Value = "QUOTIENT"
$ws.Range("E5").HorizontalAlignment = $xlRight
$ws.Range("F5").Formula = "=QUOTIENT(7, 2)"

# Row 6: ROUND
$ws.Range("E6").Value = "ROUND"
$ws.Range("E6").HorizontalAlignment = $xlRight
$ws.Range("F6").Formula = "=ROUND(18.7, 0)"

# Row 7: blank label cell
$ws.Range("E7").HorizontalAlignment = $xlRight

# Row 8: TRUE
$ws.Range("E8").Value = $true
$ws.Range("E8").HorizontalAlignment = $xlRight
$ws.Range("F8").Formula = "=TRUE()"

# Row 9: FALSE
$ws.Range("E9").Value = $false
$ws.Range("E9").HorizontalAlignment = $xlRight
$ws.Range("F9").Formula = "=FALSE()"

# Row 10: IF
$ws.Range("E10").Value = "IF"
$ws.Range("E10").HorizontalAlignment = $xlRight
$ws.Range("F10").Formula = "=IF(AbsVal > 4, 3, 1)"

# Row 11: IF ELSE
$ws.Range("E11").Value = "IF ELSE"
$ws.Range("E11").HorizontalAlignment = $xlRight
$ws.Range("F11").Formula = "=IF(AbsVal < 4, 3, 1)"

# Rows 12-22: blank label cells, right-aligned to match the column's look
foreach ($r in 12..22) {
    $ws.Range("E$r").HorizontalAlignment = $xlRight
}

# New trailing blank rows 23-24 (column E only)
$ws.Range("E23").HorizontalAlignment = $xlRight
$ws.Range("E24").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# Defined names for the new value cells
# ---------------------------------------------------------------------
$wb.Names.Add("FalseVal", "=Sheet1!`$F`$9")
$wb.Names.Add("IfElseVal", "=Sheet1!`$F`$11")
$wb.Names.Add("IfVal", "=Sheet1!`$F`$10")
$wb.Names.Add("PowOpVal", "=Sheet1!`$D`$20")
$wb.Names.Add("PowVal", "=Sheet1!`$D`$21")
$wb.Names.Add("QuoVal", "=Sheet1!`$F`$5")
$wb.Names.Add("RoundVal", "=Sheet1!`$F`$6")
$wb.Names.Add("ToDegVal", "=Sheet1!`$F`$2")
$wb.Names.Add("ToRadVal", "=Sheet1!`$F`$3")
$wb.Names.Add("TrueVal", "=Sheet1!`$F`$8")

# ---------------------------------------------------------------------
# Final selection, matching the author's saved cursor position
# ---------------------------------------------------------------------
$null = $ws.Range("F11").Select()
